# Update from Github Action
# Inserts a new row 2 (new company entry) into the survey sheet, shifting
# all existing data rows down by one. Column A (the sequential 0-based
# index) is NOT carried along with the shifted data - it is re-numbered so
# that it always equals (row number - 2), matching the original layout's
# convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a brand-new blank row at row 2; this pushes the former rows
#    2..31 down to 3..32 (dimension grows from S31 to S32).
$ws.Rows.Item(2).Insert()

# 2. The inserted row inherits formatting copied down from row 1 (bold /
#    bordered header style). Strip that back to the plain, unstyled look
#    used by every other data row.
$ws.Rows.Item(2).ClearFormats()

# 3. Column A keeps the bordered/centered numbering style; restore it on
#    the new A2 by copying the format from A3 (still carrying that style).
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# 4. Populate the new row with the survey entry for 宝马诚迈信息技术有限公司.
#    (O2/P2/R2/S2 are intentionally left blank - the freshly inserted row
#    already starts out empty there, matching the source entry.)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "宝马诚迈信息技术有限公司"
$ws.Range("C2").Value = "天隆寺金地威新"
$ws.Range("D2").Value = "开发部"
$ws.Range("E2").Value = "开发工程师"
$ws.Range("F2").Value = "9:00-18:00"
$ws.Range("G2").Value = "12:00-13:00"
$ws.Range("H2").Value = "不加班"
$ws.Range("I2").Value = "全额12%"
$ws.Range("J2").Value = "一个月"
$ws.Range("K2").Value = "3个月，不打折"
$ws.Range("L2").Value = "Macbook Pro + 4K显示器 + 升降办公桌"
$ws.Range("M2").Value = "10天起+12天病假"
$ws.Range("N2").Value = "刷脸闸机"
$ws.Range("Q2").Value = "2022-06-23 10:01:26"

# 5. Re-number column A for every shifted-down data row (now rows 3..32)
#    so it stays a plain 0-based row index: A3=1, A4=2, ... A32=30.
For ($r = 3; $r -le 32; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

Write-Host "Row insert + renumber complete"
